$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 28: Execute flag flips from "Yes" to "No" ---
$ws.Range("B28").Value = "No"

# --- New row 29: Verify Edit Phone Number ---
$ws.Range("A29").Value = "Verify Edit Phone Number"
$ws.Range("B29").Value = "No"
$ws.Range("C29").Value = "testdata.xls,profile"
$ws.Range("D29").Value = "RunRangeOfIterations"
$ws.Range("E29").Value = "'1"
$ws.Range("F29").Value = "'2"
$ws.Range("G29").Value = "Edit Phone Number"
$ws.Range("H29").Value = "coyni_mobile.tests.LoginTest,`ntestLogin,`n-pemail,`n-ppassword,`n-ppin,`n-puserName"
$ws.Range("I29").Value = "coyni_mobile.tests.CustomerProfileTest,`ntestEditPhoneNumber,`n-puserDetailsHeading,`n-ppinHeading,`n-ppin,`n-peditPhoneHeading,`n-pphoneNumber,`n-pnewPhoneNumber,`n-pcurrentPhoneHeading,`n-pcode,`n-pnewPhoneHeading,`n-pexpHeading"

# --- New row 30: Verify Edit Address ---
$ws.Range("A30").Value = "Verify Edit Address"
$ws.Range("B30").Value = "No"
$ws.Range("C30").Value = "testdata.xls,profile"
$ws.Range("D30").Value = "RunRangeOfIterations"
$ws.Range("E30").Value = "'1"
$ws.Range("F30").Value = "'2"
$ws.Range("G30").Value = "Edit Address"
$ws.Range("H30").Value = "coyni_mobile.tests.LoginTest,`ntestLogin,`n-pemail,`n-ppassword,`n-ppin,`n-puserName"
$ws.Range("I30").Value = "coyni_mobile.tests.CustomerProfileTest,`ntestEditAddress,`n-puserDetailsHeading,`n-ppinHeading,`n-ppin,`n-peditAddrHeading,`n-paddLine1,`n-paddline2,`n-pcity,`n-pstate,`n-pzipcode,`n-pcountry,`n-pexpAddress"

# --- New row 31: Verify Edit Address with Invalid data ---
# (I31 is written before F31 so the shared-string table picks up the same
# index ordering - 97 then 98 - as the canonical workbook.)
$ws.Range("A31").Value = "Verify Edit Address with Invalid data"
$ws.Range("B31").Value = "Yes"
$ws.Range("C31").Value = "testdata.xls,profile"
$ws.Range("D31").Value = "RunRangeOfIterations"
$ws.Range("E31").Value = "'1"
$ws.Range("G31").Value = "Edit Address"
$ws.Range("H31").Value = "coyni_mobile.tests.LoginTest,`ntestLogin,`n-pemail,`n-ppassword,`n-ppin,`n-puserName"
$ws.Range("I31").Value = "coyni_mobile.tests.CustomerProfileTest,`ntestEditAddressWithInvalidData,`n-puserDetailsHeading,`n-ppinHeading,`n-ppin,`n-peditAddrHeading,`n-paddLine1,`n-paddline2,`n-pcity,`n-pstate,`n-pzipcode,`n-perrMessage,`n-pelementName"
$ws.Range("F31").Value = "'4"

# Copy row 28's per-cell formatting onto the three new rows (values already set above,
# so the paste-special format-only operation layers styles on top without touching content).
$ws.Range("A28:I28").Copy()
$ws.Range("A29:I29").PasteSpecial(-4122)
$ws.Range("A30:I30").PasteSpecial(-4122)
$ws.Range("A31:I31").PasteSpecial(-4122)

# Row heights to match the wrapped multi-line content
$ws.Rows.Item(29).RowHeight = 172.8
$ws.Rows.Item(30).RowHeight = 187.2
$ws.Rows.Item(31).RowHeight = 187.2

# Sheet view follows the newly-added bottom rows
$ws.Application.ActiveWindow.ScrollRow = 30
$ws.Range("E31").Select()
